$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
# VALOR MORA total increased
$ws.Range("E11").Value = 284700
# Cant. Periodos increased from 4 to 5 (one more period/month added below)
$ws.Range("F13").Value = 5

# --- Add a new worker-period row (2509) to the table, mirroring the ---
# --- existing rows (same worker, new period), pushing the signature  ---
# --- block down by one row.                                          ---
$ws.Rows("19").Copy()
$ws.Rows("20").Insert()

# New row 20 now has the same content/format as row 19; update the period
$ws.Range("E20").Value = "2509"

$excel.CutCopyMode = 0
